# Update "想去人数" (interest count) values in column F across the sheets
# of the 广州-漫展信息 workbook, per the regenerated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 618
$ws1.Range("F5").Value  = 765
$ws1.Range("F6").Value  = 420
$ws1.Range("F10").Value = 263
$ws1.Range("F11").Value = 7020
$ws1.Range("F15").Value = 558
$ws1.Range("F16").Value = 382
$ws1.Range("F18").Value = 125
$ws1.Range("F20").Value = 728
$ws1.Range("F23").Value = 108
$ws1.Range("F25").Value = 1048
$ws1.Range("F27").Value = 41
$ws1.Range("F28").Value = 2003
$ws1.Range("F29").Value = 548
$ws1.Range("F31").Value = 542

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 278
$ws2.Range("F5").Value = 281

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 326

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 326
$ws4.Range("F3").Value  = 618
$ws4.Range("F6").Value  = 765
$ws4.Range("F8").Value  = 420
$ws4.Range("F12").Value = 263
$ws4.Range("F13").Value = 7020
$ws4.Range("F16").Value = 278
$ws4.Range("F18").Value = 558
$ws4.Range("F19").Value = 382
$ws4.Range("F22").Value = 125
$ws4.Range("F24").Value = 281
$ws4.Range("F27").Value = 728
$ws4.Range("F33").Value = 108
$ws4.Range("F35").Value = 1048
$ws4.Range("F37").Value = 41
$ws4.Range("F38").Value = 2003
$ws4.Range("F39").Value = 548
$ws4.Range("F41").Value = 542
